# Apply the "kansepur bill2" edits:
#  - Sheet1: merge duplicate name "wazid" into "wajid" (find & replace, whole-cell match)
#  - Sheet2: fill in the missing "Column sum" (Q) formulas for every worker row,
#            fix the Q7 formula (it wrongly summed Q5:Q6 instead of its own row),
#            refresh the jcb-hours figures in row 7, add a grand-total row (26)
#            labelled "Column sum" that totals the whole Q column.
#  - Leave the final selection on Sheet1 (as the workbook was last saved).

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Sheet1")
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws3 = $wb.Worksheets.Item("Sheet3")

# ---------------------------------------------------------------------------
# Sheet1: "wazid" was a duplicate spelling of "wajid" -- merge the two.
# ---------------------------------------------------------------------------
$ws1.Cells.Replace("wazid", "wajid", 2)

# ---------------------------------------------------------------------------
# Sheet2: add the missing per-row "Column sum" formulas.
# ---------------------------------------------------------------------------
$ws2.Range("Q4").Formula = "=SUM(B4:P4)"
$ws2.Range("Q7").Formula = "=SUM(B7:P7)"
$ws2.Range("Q14").Formula = "=SUM(B14:P14)"
$ws2.Range("Q17").Formula = "=SUM(B17:P17)"
$ws2.Range("Q20").Formula = "=SUM(B20:P20)"
$ws2.Range("Q23").Formula = "=SUM(B23:P23)"

# Row 7 jcb-hours figures were updated.
$ws2.Range("R7").Value = 12000
$ws2.Range("S7").Value = 35
$ws2.Range("U7").ClearContents()

# Row 20's stray R value is removed now that Q20 carries the row total.
$ws2.Range("R20").ClearContents()

# New zero placeholders next to the totals row / last data row.
$ws2.Range("R24").Value = 0
$ws2.Range("R25").Value = 0

# The grand-total formula in row 25 now also (harmlessly) includes the
# text label in column A.
$ws2.Range("Q25").Formula = "=SUM(A25:P25)"

# New row 26: label + grand total of the whole "Column sum" column.
$ws2.Range("P26").Value = "Column sum"
$ws2.Range("Q26").Formula = "=SUM(Q2:Q24)"

# ---------------------------------------------------------------------------
# Restore the view: Sheet1 active, scrolled near the bottom of the data;
# Sheet2 selection left near its newly-filled formulas.
# ---------------------------------------------------------------------------
$ws2.Activate()
$ws2.Range("Q15").Select()

$ws1.Activate()
$ws1.Range("B397").Select()
$excel.ActiveWindow.ScrollRow = 385
